# Apply updated NATMI TPM-derived values to the LR-pairs sheet (Hgf-Cd44).
# Only the data cells listed below changed between the old and new TPM runs;
# cluster labels (A:D) and unaffected numeric columns are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "0.3333333333333333"
$ws.Range("G2").Value = "0.06674466666666666"
$ws.Range("H2").Value = "0.200234"
$ws.Range("I2").Value = "0.0009912440954723497"
$ws.Range("J2").Value = "0.0009958565080158308"
$ws.Range("M2").Value = "16.14072933333334"
$ws.Range("N2").Value = "48.42218800000001"
$ws.Range("O2").Value = "0.03423048004954622"
$ws.Range("P2").Value = "0.03634868370049611"
$ws.Range("Q2").Value = "1.077307599110222"
$ws.Range("R2").Value = "9.695768391992001"
$ws.Range("S2").Value = "0.00003393076123429676"
$ws.Range("T2").Value = "0.000036198073220948"

# Row 3
$ws.Range("E3").Value = "1"
$ws.Range("F3").Value = "0.3333333333333333"
$ws.Range("G3").Value = "0.06674466666666666"
$ws.Range("H3").Value = "0.200234"
$ws.Range("I3").Value = "0.0009912440954723497"
$ws.Range("J3").Value = "0.0009958565080158308"
$ws.Range("O3").Value = "0.1719151703242873"
$ws.Range("P3").Value = "0.1825533892714798"
$ws.Range("Q3").Value = "5.410544027562889"
$ws.Range("R3").Value = "48.694896248066"
$ws.Range("S3").Value = "0.0001704098975060731"
$ws.Range("T3").Value = "0.0001817969807663506"

# Row 4
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = "0.3333333333333333"
$ws.Range("G4").Value = "0.06674466666666666"
$ws.Range("H4").Value = "0.200234"
$ws.Range("I4").Value = "0.0009912440954723497"
$ws.Range("J4").Value = "0.0009958565080158308"
$ws.Range("M4").Value = "168.70371"
$ws.Range("N4").Value = "506.11113"
$ws.Range("O4").Value = "0.3577786889414888"
$ws.Range("P4").Value = "0.3799182594076638"
$ws.Range("Q4").Value = "11.26007288938"
$ws.Range("R4").Value = "101.34065600442"
$ws.Range("S4").Value = "0.0003546460128990892"
$ws.Range("T4").Value = "0.0003783440711451687"

# Row 5
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "0.3333333333333333"
$ws.Range("G5").Value = "0.06674466666666666"
$ws.Range("H5").Value = "0.200234"
$ws.Range("I5").Value = "0.0009912440954723497"
$ws.Range("J5").Value = "0.0009958565080158308"
$ws.Range("M5").Value = "82.43477250000001"
$ws.Range("N5").Value = "164.869545"
$ws.Range("O5").Value = "0.1748236883957081"
$ws.Range("P5").Value = "0.1237612588479007"
$ws.Range("Q5").Value = "5.502081412255"
$ws.Range("R5").Value = "33.01248847353001"
$ws.Range("S5").Value = "0.0001732929488709436"
$ws.Range("T5").Value = "0.0001232484550639137"

# Row 6
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "0.3333333333333333"
$ws.Range("G6").Value = "0.06674466666666666"
$ws.Range("H6").Value = "0.200234"
$ws.Range("I6").Value = "0.0009912440954723497"
$ws.Range("J6").Value = "0.0009958565080158308"
$ws.Range("M6").Value = "123.1883796666667"
$ws.Range("N6").Value = "369.565139"
$ws.Range("O6").Value = "0.2612519722889696"
$ws.Range("P6").Value = "0.2774184087724594"
$ws.Range("Q6").Value = "8.222167338058444"
$ws.Range("R6").Value = "73.99950604252599"
$ws.Range("S6").Value = "0.000258964474961947"
$ws.Range("T6").Value = "0.0002762689278194498"

# Row 7
$ws.Range("I7").Value = "0.1187608236941705"
$ws.Range("J7").Value = "0.1193134362296531"
$ws.Range("M7").Value = "16.14072933333334"
$ws.Range("N7").Value = "48.42218800000001"
$ws.Range("O7").Value = "0.03423048004954622"
$ws.Range("P7").Value = "0.03634868370049611"
$ws.Range("Q7").Value = "129.0720806577436"
$ws.Range("R7").Value = "1161.648725919692"
$ws.Range("S7").Value = "0.004065240006130979"
$ws.Range("T7").Value = "0.004336886354730973"

# Row 8
$ws.Range("I8").Value = "0.1187608236941705"
$ws.Range("J8").Value = "0.1193134362296531"
$ws.Range("O8").Value = "0.1719151703242873"
$ws.Range("P8").Value = "0.1825533892714798"
$ws.Range("S8").Value = "0.02041678723323597"
$ws.Range("T8").Value = "0.02178107216934974"

# Row 9
$ws.Range("I9").Value = "0.1187608236941705"
$ws.Range("J9").Value = "0.1193134362296531"
$ws.Range("M9").Value = "168.70371"
$ws.Range("N9").Value = "506.11113"
$ws.Range("O9").Value = "0.3577786889414888"
$ws.Range("P9").Value = "0.3799182594076638"
$ws.Range("Q9").Value = "1349.06784041113"
$ws.Range("R9").Value = "12141.61056370017"
$ws.Range("S9").Value = "0.04249009179891161"
$ws.Range("T9").Value = "0.04532935301631709"

# Row 10
$ws.Range("I10").Value = "0.1187608236941705"
$ws.Range("J10").Value = "0.1193134362296531"
$ws.Range("M10").Value = "82.43477250000001"
$ws.Range("N10").Value = "164.869545"
$ws.Range("O10").Value = "0.1748236883957081"
$ws.Range("P10").Value = "0.1237612588479007"
$ws.Range("Q10").Value = "659.2036447293176"
$ws.Range("R10").Value = "3955.221868375906"
$ws.Range("S10").Value = "0.02076220523512728"
$ws.Range("T10").Value = "0.01476638106525059"

# Row 11
$ws.Range("I11").Value = "0.1187608236941705"
$ws.Range("J11").Value = "0.1193134362296531"
$ws.Range("M11").Value = "123.1883796666667"
$ws.Range("N11").Value = "369.565139"
$ws.Range("O11").Value = "0.2612519722889696"
$ws.Range("P11").Value = "0.2774184087724594"
$ws.Range("Q11").Value = "985.0967789662501"
$ws.Range("R11").Value = "8865.871010696252"
$ws.Range("S11").Value = "0.03102649942076462"
$ws.Range("T11").Value = "0.03309974362400467"

# Row 12
$ws.Range("G12").Value = "24.06383433333333"
$ws.Range("H12").Value = "72.191503"
$ws.Range("I12").Value = "0.3573788721796719"
$ws.Range("J12").Value = "0.3590418115105046"
$ws.Range("M12").Value = "16.14072933333334"
$ws.Range("N12").Value = "48.42218800000001"
$ws.Range("O12").Value = "0.03423048004954622"
$ws.Range("P12").Value = "0.03634868370049611"
$ws.Range("Q12").Value = "388.4078366965072"
$ws.Range("R12").Value = "3495.670530268564"
$ws.Range("S12").Value = "0.01223325035427559"
$ws.Range("T12").Value = "0.01305069724184848"

# Row 13
$ws.Range("G13").Value = "24.06383433333333"
$ws.Range("H13").Value = "72.191503"
$ws.Range("I13").Value = "0.3573788721796719"
$ws.Range("J13").Value = "0.3590418115105046"
$ws.Range("O13").Value = "0.1719151703242873"
$ws.Range("P13").Value = "0.1825533892714798"
$ws.Range("Q13").Value = "1950.694214755927"
$ws.Range("R13").Value = "17556.24793280335"
$ws.Range("S13").Value = "0.06143884968107"
$ws.Range("T13").Value = "0.06554429958141443"

# Row 14
$ws.Range("G14").Value = "24.06383433333333"
$ws.Range("H14").Value = "72.191503"
$ws.Range("I14").Value = "0.3573788721796719"
$ws.Range("J14").Value = "0.3590418115105046"
$ws.Range("M14").Value = "168.70371"
$ws.Range("N14").Value = "506.11113"
$ws.Range("O14").Value = "0.3577786889414888"
$ws.Range("P14").Value = "0.3799182594076638"
$ws.Range("Q14").Value = "4059.65812885871"
$ws.Range("R14").Value = "36536.92315972839"
$ws.Range("S14").Value = "0.1278625443438309"
$ws.Range("T14").Value = "0.1364065400836454"

# Row 15
$ws.Range("G15").Value = "24.06383433333333"
$ws.Range("H15").Value = "72.191503"
$ws.Range("I15").Value = "0.3573788721796719"
$ws.Range("J15").Value = "0.3590418115105046"
$ws.Range("M15").Value = "82.43477250000001"
$ws.Range("N15").Value = "164.869545"
$ws.Range("O15").Value = "0.1748236883957081"
$ws.Range("P15").Value = "0.1237612588479007"
$ws.Range("Q15").Value = "1983.696708746023"
$ws.Range("R15").Value = "11902.18025247614"
$ws.Range("S15").Value = "0.06247829258914854"
$ws.Range("T15").Value = "0.04443546657157073"

# Row 16
$ws.Range("G16").Value = "24.06383433333333"
$ws.Range("H16").Value = "72.191503"
$ws.Range("I16").Value = "0.3573788721796719"
$ws.Range("J16").Value = "0.3590418115105046"
$ws.Range("M16").Value = "123.1883796666667"
$ws.Range("N16").Value = "369.565139"
$ws.Range("O16").Value = "0.2612519722889696"
$ws.Range("P16").Value = "0.2774184087724594"
$ws.Range("Q16").Value = "2964.384760090435"
$ws.Range("R16").Value = "26679.46284081391"
$ws.Range("S16").Value = "0.09336593521134684"
$ws.Range("T16").Value = "0.09960480803202548"

# Row 17
$ws.Range("G17").Value = "0.9355965000000001"
$ws.Range("H17").Value = "1.871193"
$ws.Range("I17").Value = "0.01389481066706348"
$ws.Range("J17").Value = "0.009306310251024633"
$ws.Range("M17").Value = "16.14072933333334"
$ws.Range("N17").Value = "48.42218800000001"
$ws.Range("O17").Value = "0.03423048004954622"
$ws.Range("P17").Value = "0.03634868370049611"
$ws.Range("Q17").Value = "15.101209871714"
$ws.Range("R17").Value = "90.60725923028401"
$ws.Range("S17").Value = "0.0004756260393311385"
$ws.Range("T17").Value = "0.000338272127733179"

# Row 18
$ws.Range("G18").Value = "0.9355965000000001"
$ws.Range("H18").Value = "1.871193"
$ws.Range("I18").Value = "0.01389481066706348"
$ws.Range("J18").Value = "0.009306310251024633"
$ws.Range("O18").Value = "0.1719151703242873"
$ws.Range("P18").Value = "0.1825533892714798"
$ws.Range("Q18").Value = "75.8425550398595"
$ws.Range("R18").Value = "455.0553302391571"
$ws.Range("S18").Value = "0.002388728742451943"
$ws.Range("T18").Value = "0.001698898477936463"

# Row 19
$ws.Range("G19").Value = "0.9355965000000001"
$ws.Range("H19").Value = "1.871193"
$ws.Range("I19").Value = "0.01389481066706348"
$ws.Range("J19").Value = "0.009306310251024633"
$ws.Range("M19").Value = "168.70371"
$ws.Range("N19").Value = "506.11113"
$ws.Range("O19").Value = "0.3577786889414888"
$ws.Range("P19").Value = "0.3799182594076638"
$ws.Range("Q19").Value = "157.838600613015"
$ws.Range("R19").Value = "947.0316036780901"
$ws.Range("S19").Value = "0.004971267143552185"
$ws.Range("T19").Value = "0.003535637192076978"

# Row 20
$ws.Range("G20").Value = "0.9355965000000001"
$ws.Range("H20").Value = "1.871193"
$ws.Range("I20").Value = "0.01389481066706348"
$ws.Range("J20").Value = "0.009306310251024633"
$ws.Range("M20").Value = "82.43477250000001"
$ws.Range("N20").Value = "164.869545"
$ws.Range("O20").Value = "0.1748236883957081"
$ws.Range("P20").Value = "0.1237612588479007"
$ws.Range("Q20").Value = "77.12568462929626"
$ws.Range("R20").Value = "308.502738517185"
$ws.Range("S20").Value = "0.002429142050376067"
$ws.Range("T20").Value = "0.001151760671895931"

# Row 21
$ws.Range("G21").Value = "0.9355965000000001"
$ws.Range("H21").Value = "1.871193"
$ws.Range("I21").Value = "0.01389481066706348"
$ws.Range("J21").Value = "0.009306310251024633"
$ws.Range("M21").Value = "123.1883796666667"
$ws.Range("N21").Value = "369.565139"
$ws.Range("O21").Value = "0.2612519722889696"
$ws.Range("P21").Value = "0.2774184087724594"
$ws.Range("Q21").Value = "115.2546168568045"
$ws.Range("R21").Value = "691.527701140827"
$ws.Range("S21").Value = "0.003630046691352147"
$ws.Range("T21").Value = "0.002581741781382081"

# Row 22
$ws.Range("G22").Value = "34.27139366666666"
$ws.Range("H22").Value = "102.814181"
$ws.Range("I22").Value = "0.5089742493636218"
$ws.Range("J22").Value = "0.5113425855008019"
$ws.Range("M22").Value = "16.14072933333334"
$ws.Range("N22").Value = "48.42218800000001"
$ws.Range("O22").Value = "0.03423048004954622"
$ws.Range("P22").Value = "0.03634868370049611"
$ws.Range("Q22").Value = "553.1652890497809"
$ws.Range("R22").Value = "4978.487601448028"
$ws.Range("S22").Value = "0.01742243288857422"
$ws.Range("T22").Value = "0.01858662990296254"

# Row 23
$ws.Range("G23").Value = "34.27139366666666"
$ws.Range("H23").Value = "102.814181"
$ws.Range("I23").Value = "0.5089742493636218"
$ws.Range("J23").Value = "0.5113425855008019"
$ws.Range("O23").Value = "0.1719151703242873"
$ws.Range("P23").Value = "0.1825533892714798"
$ws.Range("Q23").Value = "2778.152825985196"
$ws.Range("R23").Value = "25003.37543386677"
$ws.Range("S23").Value = "0.08750039477002333"
$ws.Range("T23").Value = "0.09334732206201286"

# Row 24
$ws.Range("G24").Value = "34.27139366666666"
$ws.Range("H24").Value = "102.814181"
$ws.Range("I24").Value = "0.5089742493636218"
$ws.Range("J24").Value = "0.5113425855008019"
$ws.Range("M24").Value = "168.70371"
$ws.Range("N24").Value = "506.11113"
$ws.Range("O24").Value = "0.3577786889414888"
$ws.Range("P24").Value = "0.3799182594076638"
$ws.Range("Q24").Value = "5781.711258437169"
$ws.Range("R24").Value = "52035.40132593452"
$ws.Range("S24").Value = "0.182100139642295"
$ws.Range("T24").Value = "0.1942683850444792"

# Row 25
$ws.Range("G25").Value = "34.27139366666666"
$ws.Range("H25").Value = "102.814181"
$ws.Range("I25").Value = "0.5089742493636218"
$ws.Range("J25").Value = "0.5113425855008019"
$ws.Range("M25").Value = "82.43477250000001"
$ws.Range("N25").Value = "164.869545"
$ws.Range("O25").Value = "0.1748236883957081"
$ws.Range("P25").Value = "0.1237612588479007"
$ws.Range("Q25").Value = "2825.154540169607"
$ws.Range("R25").Value = "16950.92724101765"
$ws.Range("S25").Value = "0.08898075557218522"
$ws.Range("T25").Value = "0.06328440208411953"

# Row 26
$ws.Range("G26").Value = "34.27139366666666"
$ws.Range("H26").Value = "102.814181"
$ws.Range("I26").Value = "0.5089742493636218"
$ws.Range("J26").Value = "0.5113425855008019"
$ws.Range("M26").Value = "123.1883796666667"
$ws.Range("N26").Value = "369.565139"
$ws.Range("O26").Value = "0.2612519722889696"
$ws.Range("P26").Value = "0.2774184087724594"
$ws.Range("Q26").Value = "4221.837454715128"
$ws.Range("R26").Value = "37996.53709243616"
$ws.Range("S26").Value = "0.132970526490544"
$ws.Range("T26").Value = "0.1418558464072278"
